$d = $word.ActiveDocument

# --- Change 1: "This application returns current weather data for any
#     location specified by the user" -> append a trailing period. ---
$d.Content.Find.Execute(
    "This application returns current weather data for any location specified by the user",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This application returns current weather data for any location specified by the user.",
    2) | Out-Null

# --- Change 2: "amount of CO2 emissions for" -> "amount of CO2 emissions from"
#     and append a trailing period to the CO2 paragraph. ---
$d.Content.Find.Execute(
    "amount of CO2 emissions for new motor vehicles to ensure government compliance",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "amount of CO2 emissions from new motor vehicles to ensure government compliance.",
    2) | Out-Null

# --- Change 3: the empty bullet paragraph right after the CO2-emissions
#     bullet (currently blank, ilvl 0) becomes the "Logistic Model for
#     Predicting Late E-Commerce Deliveries" heading, followed by a new
#     ilvl-1 description paragraph. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains("new motor vehicles to ensure government compliance")) {
        $emptyPara = $d.Paragraphs($i + 1)
        $emptyPara.Range.Text = "Logistic Model for Predicting Late E-Commerce Deliveries "

        $emptyPara.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 2)
        $newPara.Range.ListFormat.ListLevelNumber = 2
        $newPara.Range.Text = "This model predicts which deliveries from an online retailer will be late so the retailer can remedy the problem accordingly."
        break
    }
}

# --- Change 4: after "...the IMDb Rating of New Movies" add a new ilvl-1
#     description paragraph. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains("the IMDb Rating of New Movies")) {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.ListFormat.ListLevelNumber = 2
        $newPara.Range.Text = "This model predicts the IMDb rating for new films allowing directors to focus on films predicted to be good ones."
        break
    }
}

Write-Host "Done."
